# Region I_ELECTRIFICATION.xlsx update
# - Inserts two new columns at AS (shifting the old "Status as of July 4, 2025"
#   column from AS to AU), giving two fresh blank columns AS ("Unnamed: 44")
#   and AT ("Unnamed: 45").
# - Moves the existing AR2:AR303 status values (COMPLETED/REVERTED/ONGOING)
#   down into the new AT column.
# - For rows whose project is an on-grid energization, stamps the now-empty
#   AS cell with "ongrid".
# - For rows that belong to the BBM batches, stamps the now-empty AR cell
#   with the appropriate BBM label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two blank columns at column 45 (AS). Excel shifts the previous
#    AS column (and everything after it) two places to the right, so the old
#    AS1 header ("Status as of July 4, 2025") ends up in AU1, and the
#    dimension / dataValidation ranges are auto-adjusted by the engine.
$ws.Columns.Item(45).Insert()
$ws.Columns.Item(45).Insert()

# 2) New header labels for the freshly inserted columns (they inherit AR's
#    bold/bordered style automatically from the Insert above).
$ws.Range("AS1").Value = "Unnamed: 44"
$ws.Range("AT1").Value = "Unnamed: 45"

# 3) Move the old AR status values (rows 2-303) over to the new AT column.
$ws.Range("AR2:AR303").Cut($ws.Range("AT2:AT303"))

# 4) Rows whose AS cell should read "ongrid".
$ongridRows = @(2,3,4,5,6,7,8,18,20,21,55,56,57,58,59,60,61,62,63,64,65,66,69,70,71,72,73,75,76,77,78,79,80,81,82,83,84,85,86,88,89,90,91,92,93,97,100,101,106,110,112,114,127,129)
foreach ($r in $ongridRows) {
    $ws.Range("AS$r").Value = "ongrid"
}

# 5) Rows whose (now empty) AR cell should carry a BBM batch label.
$bbmLabels = @{
    223 = "BBM 2025 UPGRADE"
    227 = "BBM 2025 UPGRADE"
    228 = "BBM 2025 UPGRADE"
    229 = "BBM 2025 UPGRADE"
    238 = "BBM 2025 UPGRADE"
    241 = "BBM 2023 UPGRADE"
    242 = "BBM 2023 UPGRADE"
    246 = "BBM 2023 UPGRADE"
    250 = "BBM 2023 UPGRADE"
    254 = "bbm 2023 ONGRID"
    261 = "bbm 2023 ONGRID"
    264 = "BBM 2023 UPGRADE"
    267 = "BBM 2023 UPGRADE"
    269 = "BBM 2023 UPGRADE"
    271 = "bbm 2023 ONGRID"
    279 = "BBM 2025 UPGRADE"
    286 = "BBM 2023 UPGRADE"
    298 = "BBM 2024 SOLAR"
    300 = "BBM 2025 ONGRID"
    303 = "BBM 2024 ONGRID"
}
foreach ($r in $bbmLabels.Keys) {
    $ws.Range("AR$r").Value = $bbmLabels[$r]
}
